# Apply updated harmonic similarity match data (rows 2-17) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'schubert-winterreise_133'
$ws.Cells.Item(2, 2).Value = 'isophonics_47'
$ws.Cells.Item(2, 3).Value = 0.2559523809523809
$ws.Cells.Item(2, 4).Value = '[[''G:maj/B'', ''C:maj'', ''G:maj/D'']]'
$ws.Cells.Item(2, 5).Value = '[[''A'', ''D'', ''A'']]'
$ws.Cells.Item(2, 6).Value = '[(61.6, 63.84)]'
$ws.Cells.Item(2, 7).Value = '[(1.302038, 5.439841)]'

# Row 3
$ws.Cells.Item(3, 1).Value = 'jaah_9'
$ws.Cells.Item(3, 2).Value = 'jaah_49'
$ws.Cells.Item(3, 3).Value = 0.1125385405960945
$ws.Cells.Item(3, 4).Value = '[[''C:7'', ''F'', ''F''], [''G:min7'', ''C:7'', ''F:7'']]'
$ws.Cells.Item(3, 5).Value = '[[''F:7'', ''Bb'', ''Bb''], [''C:min7'', ''F:7'', ''Bb:7'']]'
$ws.Cells.Item(3, 6).Value = '[(51.69, 55.72), (3.05, 6.96)]'
$ws.Cells.Item(3, 7).Value = '[(44.74, 47.14), (2.6, 3.8)]'
$ws.Cells.Item(3, 9).ClearContents()

# Row 4
$ws.Cells.Item(4, 1).Value = 'isophonics_115'
$ws.Cells.Item(4, 2).Value = 'isophonics_49'
$ws.Cells.Item(4, 3).Value = 0.1047619047619048
$ws.Cells.Item(4, 4).Value = '[[''Bb'', ''Bb:min'', ''F'']]'
$ws.Cells.Item(4, 5).Value = '[[''C'', ''C:min'', ''G'']]'
$ws.Cells.Item(4, 6).Value = '[(39.61044, 45.427039)]'
$ws.Cells.Item(4, 7).Value = '[(24.911369, 30.797628)]'
$ws.Cells.Item(4, 8).ClearContents()

# Row 5
$ws.Cells.Item(5, 1).Value = 'schubert-winterreise_49'
$ws.Cells.Item(5, 2).Value = 'schubert-winterreise_177'
$ws.Cells.Item(5, 3).Value = 0.2913752913752914
$ws.Cells.Item(5, 4).Value = '[[''A:min'', ''E:maj/G#'', ''A:min'', ''A:7/G'', ''D:min/F'']]'
$ws.Cells.Item(5, 5).Value = '[[''A:min'', ''E:maj'', ''A:min'', ''A:7'', ''D:min/A'']]'
$ws.Cells.Item(5, 6).Value = '[(13.34, 17.36)]'
$ws.Cells.Item(5, 7).Value = '[(13.5, 25.26)]'
$ws.Cells.Item(5, 9).Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'

# Row 6
$ws.Cells.Item(6, 1).Value = 'isophonics_226'
$ws.Cells.Item(6, 2).Value = 'isophonics_19'
$ws.Cells.Item(6, 3).Value = 0.225
$ws.Cells.Item(6, 4).Value = '[[''F#:min'', ''B'', ''E'']]'
$ws.Cells.Item(6, 5).Value = '[[''G:min'', ''C'', ''F'']]'
$ws.Cells.Item(6, 6).Value = '[(45.645192, 58.300068)]'
$ws.Cells.Item(6, 7).Value = '[(45.474603, 48.945986)]'

# Row 7
$ws.Cells.Item(7, 1).Value = 'schubert-winterreise_113'
$ws.Cells.Item(7, 2).Value = 'schubert-winterreise_4'
$ws.Cells.Item(7, 3).Value = 0.1833333333333333
$ws.Cells.Item(7, 4).Value = '[[''C:7'', ''F:min'', ''C:maj''], [''F:maj'', ''A#:maj'', ''F:maj''], [''F:min'', ''C:maj'', ''F:min'']]'
$ws.Cells.Item(7, 5).Value = '[[''A:7'', ''D:min/A'', ''A:maj''], [''D:maj'', ''G:maj'', ''D:maj/F#''], [''D:min'', ''A:maj'', ''D:min'']]'
$ws.Cells.Item(7, 6).Value = '[(12.96, 17.7), (123.22, 126.88), (0.78, 5.26)]'
$ws.Cells.Item(7, 7).Value = '[(31.98, 34.58), (66.52, 70.14), (8.24, 13.46)]'
$ws.Cells.Item(7, 8).Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'

# Row 8
$ws.Cells.Item(8, 1).Value = 'schubert-winterreise_53'
$ws.Cells.Item(8, 2).Value = 'jaah_25'
$ws.Cells.Item(8, 3).Value = 0.2729885057471264
$ws.Cells.Item(8, 4).Value = '[[''E:maj'', ''B:7/A'', ''E:maj/G#'', ''E:maj/B''], [''B:7'', ''E:maj'', ''B:7'', ''E:maj'']]'
$ws.Cells.Item(8, 5).Value = '[[''Bb'', ''F:7'', ''Bb'', ''Bb''], [''F:7'', ''Bb'', ''F:7'', ''Bb'']]'
$ws.Cells.Item(8, 6).Value = '[(215.06, 226.54), (15.26, 37.4)]'
$ws.Cells.Item(8, 7).Value = '[(50.13, 55.78), (48.25, 53.91)]'

# Row 9
$ws.Cells.Item(9, 1).Value = 'schubert-winterreise_97'
$ws.Cells.Item(9, 2).Value = 'schubert-winterreise_92'
$ws.Cells.Item(9, 3).Value = 0.2196969696969697
$ws.Cells.Item(9, 4).Value = '[[''D:min'', ''A:maj'', ''D:min'']]'
$ws.Cells.Item(9, 5).Value = '[[''B:min'', ''F#:maj/A#'', ''B:min'']]'
$ws.Cells.Item(9, 6).Value = '[(8.54, 14.08)]'
$ws.Cells.Item(9, 7).Value = '[(13.56, 16.28)]'

# Row 10
$ws.Cells.Item(10, 1).Value = 'schubert-winterreise_111'
$ws.Cells.Item(10, 2).Value = 'schubert-winterreise_2'
$ws.Cells.Item(10, 3).Value = 0.2657342657342657
$ws.Cells.Item(10, 4).Value = '[[''G:maj'', ''D:7/C'', ''G:maj/B'', ''D:7/C'', ''G:maj/B'']]'
$ws.Cells.Item(10, 5).Value = '[[''A:maj/E'', ''E:7'', ''A:maj'', ''E:7'', ''A:maj'']]'
$ws.Cells.Item(10, 6).Value = '[(67.76, 84.78)]'
$ws.Cells.Item(10, 7).Value = '[(20.56, 26.4)]'
$ws.Cells.Item(10, 8).Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
$ws.Cells.Item(10, 9).Value = 'spotify:track:0XfunCHFEeQnzm4NaY8rJr'

# Row 11
$ws.Cells.Item(11, 1).Value = 'schubert-winterreise_128'
$ws.Cells.Item(11, 2).Value = 'schubert-winterreise_134'
$ws.Cells.Item(11, 3).Value = 0.323076923076923
$ws.Cells.Item(11, 4).Value = '[[''G:maj'', ''C:maj/G'', ''G:maj'']]'
$ws.Cells.Item(11, 5).Value = '[[''C:maj/G'', ''F:maj'', ''C:maj/G'']]'
$ws.Cells.Item(11, 6).Value = '[(18.32, 25.82)]'
$ws.Cells.Item(11, 7).Value = '[(142.14, 146.8)]'
$ws.Cells.Item(11, 8).Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'

# Row 12
$ws.Cells.Item(12, 1).Value = 'schubert-winterreise_59'
$ws.Cells.Item(12, 2).Value = 'isophonics_140'
$ws.Cells.Item(12, 3).Value = 0.2708333333333333
$ws.Cells.Item(12, 4).Value = '[[''A:maj'', ''E:maj'', ''B:maj''], [''E:maj/G#'', ''A:maj'', ''E:maj'']]'
$ws.Cells.Item(12, 5).Value = '[[''C'', ''G'', ''D''], [''G'', ''C'', ''G'']]'
$ws.Cells.Item(12, 6).Value = '[(63.08, 71.12), (59.22, 66.88)]'
$ws.Cells.Item(12, 7).Value = '[(59.224321, 65.737518), (0.465952, 5.272756)]'
$ws.Cells.Item(12, 9).Value = 'spotify:track:3VbGCXWRiouAq8VyMYN2MI'

# Row 13
$ws.Cells.Item(13, 1).Value = 'schubert-winterreise_95'
$ws.Cells.Item(13, 2).Value = 'schubert-winterreise_171'
$ws.Cells.Item(13, 3).Value = 0.1202898550724638
$ws.Cells.Item(13, 4).Value = '[[''C:min/G'', ''G:7'', ''C:min'', ''C:min''], [''C'', ''C/G'', ''G:7'', ''C'']]'
$ws.Cells.Item(13, 5).Value = '[[''F#:min'', ''C#:7'', ''F#:min'', ''F#:min''], [''F#:maj'', ''F#:maj/A#'', ''C#:7'', ''F#:maj'']]'
$ws.Cells.Item(13, 6).Value = '[(20.92, 26.04), (225.32, 230.46)]'
$ws.Cells.Item(13, 7).Value = '[(3.82, 8.62), (21.02, 24.16)]'
$ws.Cells.Item(13, 8).ClearContents()
$ws.Cells.Item(13, 9).Value = 'spotify:track:4lrfYSnZmpXdCWuWqVo8L0'

# Row 14
$ws.Cells.Item(14, 1).Value = 'isophonics_82'
$ws.Cells.Item(14, 2).Value = 'isophonics_2'
$ws.Cells.Item(14, 3).Value = 0.09441489361702127
$ws.Cells.Item(14, 4).Value = '[[''F#:min'', ''C#:min'', ''F#:min''], [''A'', ''D'', ''A''], [''E'', ''A'', ''D'']]'
$ws.Cells.Item(14, 5).Value = '[[''C:min'', ''G:min'', ''C:min''], [''Bb/3'', ''Eb:maj'', ''Bb/3''], [''Bb:maj'', ''Eb:maj'', ''Ab:maj'']]'
$ws.Cells.Item(14, 6).Value = '[(17.839297, 23.087006), (3.988594, 7.754783), (54.526825, 61.516031)]'
$ws.Cells.Item(14, 7).Value = '[(312.842, 317.997), (38.041, 42.375), (237.731, 240.233)]'
$ws.Cells.Item(14, 8).Value = 'spotify:track:5EzvwjFwdP5Kfl5AZAemzu'

# Row 15
$ws.Cells.Item(15, 1).Value = 'schubert-winterreise_185'
$ws.Cells.Item(15, 2).Value = 'schubert-winterreise_143'
$ws.Cells.Item(15, 3).Value = 1
$ws.Cells.Item(15, 4).Value = '[[''B:min'', ''B:7'', ''E:maj/B'', ''E:min/B'', ''B:maj'', ''B:min'', ''F#:maj'', ''B:min'', ''B:7'', ''E:min/B'', ''F#:7/B'', ''B:min'', ''E:min/B'', ''B:min'', ''D:maj/A'', ''E:min/G'', ''D:maj/F#'', ''A:7'', ''D:maj'', ''B:min'', ''B:7/A'', ''E:maj/G#'', ''A:(3,5,b7,b9)/G'', ''D:maj/F#'', ''F#:(3,5,b7,b9)/E'', ''G:(3,5)'', ''B:min/F#'', ''F#:7'', ''G:(3,5)'', ''B:min/F#'', ''F#:7'', ''B:min'', ''B:7'', ''E:maj'', ''E:min'', ''B:maj'', ''B:min/B'', ''F#:maj'', ''B:min'']]'
$ws.Cells.Item(15, 5).Value = '[[''B:min'', ''B:7'', ''E:maj/B'', ''E:min/B'', ''B:maj'', ''B:min'', ''F#:maj'', ''B:min'', ''B:7'', ''E:min/B'', ''F#:7/B'', ''B:min'', ''E:min/B'', ''B:min'', ''D:maj/A'', ''E:min/G'', ''D:maj/F#'', ''A:7'', ''D:maj'', ''B:min'', ''B:7/A'', ''E:maj/G#'', ''A:(3,5,b7,b9)/G'', ''D:maj/F#'', ''F#:(3,5,b7,b9)/E'', ''G:(3,5)'', ''B:min/F#'', ''F#:7'', ''G:(3,5)'', ''B:min/F#'', ''F#:7'', ''B:min'', ''B:7'', ''E:maj'', ''E:min'', ''B:maj'', ''B:min/B'', ''F#:maj'', ''B:min'']]'
$ws.Cells.Item(15, 6).Value = '[(0.82, 107.92)]'
$ws.Cells.Item(15, 7).Value = '[(1.66, 97.0)]'
$ws.Cells.Item(15, 8).Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'
$ws.Cells.Item(15, 9).Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'

# Row 16
$ws.Cells.Item(16, 1).Value = 'schubert-winterreise_48'
$ws.Cells.Item(16, 2).Value = 'isophonics_45'
$ws.Cells.Item(16, 3).Value = 0.1396103896103896
$ws.Cells.Item(16, 4).Value = '[[''F:maj'', ''F:7'', ''A#:maj'']]'
$ws.Cells.Item(16, 5).Value = '[[''F:maj'', ''F:7'', ''Bb'']]'
$ws.Cells.Item(16, 6).Value = '[(14.48, 22.82)]'
$ws.Cells.Item(16, 7).Value = '[(13.155, 20.379)]'
$ws.Cells.Item(16, 8).ClearContents()

# Row 17
$ws.Cells.Item(17, 1).Value = 'schubert-winterreise_149'
$ws.Cells.Item(17, 2).Value = 'schubert-winterreise_53'
$ws.Cells.Item(17, 3).Value = 0.2223837209302326
$ws.Cells.Item(17, 4).Value = '[[''A:min/E'', ''E:(3,5,b7,b9)'', ''A:min'']]'
$ws.Cells.Item(17, 5).Value = '[[''E:min'', ''B:(3,5,b7,b9)'', ''E:min'']]'
$ws.Cells.Item(17, 6).Value = '[(63.12, 69.42)]'
$ws.Cells.Item(17, 7).Value = '[(94.78, 114.14)]'
$ws.Cells.Item(17, 8).Value = 'spotify:track:2qCvEz2hEb92VFATqVvrht'
$ws.Cells.Item(17, 9).ClearContents()
